$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 9 (data set now only has 7 data rows, rows 2-8)
$ws.Rows.Item(9).Delete()

# Update column widths.
# Note: Excel's ColumnWidth (character units) round-trips to the stored
# OOXML <col width="..."> value with a constant +5/6 padding offset in
# this engine, so subtract it here to land exactly on the target widths
# (36, 84, 16, 25) once saved.
$ws.Columns.Item(3).ColumnWidth = 36 - (5/6)
$ws.Columns.Item(4).ColumnWidth = 84 - (5/6)
$ws.Columns.Item(6).ColumnWidth = 16 - (5/6)
$ws.Columns.Item(8).ColumnWidth = 25 - (5/6)

# New data for rows 2-8
$data = @(
    @("1331697", "https://aiesec.org/opportunity/global-talent/1331697", "Social Media Manager", "Sfax, Tunisie", "No", "0 applicants", "9 - 12 Weeks", "L'adé-s"),
    @("1331651", "https://aiesec.org/opportunity/global-talent/1331651", "Web Developer", "Sfax, Tunisie", "No", "0 applicants", "9 - 12 Weeks", "L'adé-s"),
    @("1331631", "https://aiesec.org/opportunity/global-talent/1331631", "developer mobile", "Sfax, Tunisie", "No", "2 applicants", "9 - 12 Weeks", "Sky Academy"),
    @("1331470", "https://aiesec.org/opportunity/global-talent/1331470", "CONSTRUCTION PROJECT COORDINATOR", "Abidjan, Côte d'Ivoire", "No", "1 applicant", "3 - 6 Months", "ONG REFUGE DES ENFANTS"),
    @("1328962", "https://aiesec.org/opportunity/global-talent/1328962", "Social Media Manager", "Hong Kong", "No", "46 applicants", "6 - 18 Months", "Wong's Limited"),
    @("1328685", "https://aiesec.org/opportunity/global-talent/1328685", "Medical Advisor (Russian Speaker)", "İstanbul, Türkiye", "No", "7 applicants", "6 - 18 Months", "International Plus"),
    @("1327889", "https://aiesec.org/opportunity/global-talent/1327889", "Graphic Designer", "Birkat as SAB, Madinet Berkat as Sabee, Birket el Sab, Menofia Governorate, Egypt", "No", "4 applicants", "9 - 12 Weeks", "Lines")
)

$rowIndex = 2
foreach ($rowValues in $data) {
    $colIndex = 1
    foreach ($val in $rowValues) {
        $cell = $ws.Cells.Item($rowIndex, $colIndex)
        if ($colIndex -eq 1) {
            # Column A (opportunity id) holds numeric-looking text in the
            # source data (e.g. "1331697"). Assigning it directly via
            # .Value would let Excel auto-coerce it to a number, so force
            # a text format first, then clear the formatting delta back
            # off so no stray style is left applied to the cell.
            $cell.NumberFormat = "@"
            $cell.Value = $val
            $cell.ClearFormats()
        } else {
            $cell.Value = $val
        }
        $colIndex++
    }
    $rowIndex++
}
